# Refresh cached Moogle/Universalis market-board price snapshots and their
# dependent Leve profit columns (H:N) across the crafting-job sheets.
# Generated from the authoritative cell-level diff; one block per sheet/row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 3056.4666
$ws.Range("I70").Value = 4050
$ws.Range("J70").Value = 2695.182
$ws.Range("K70").Value = 12150
$ws.Range("L70").Value = 8085.545999999999
$ws.Range("M70").Value = -11880
$ws.Range("N70").Value = -8625.545999999998
# Row 73
$ws.Range("H73").Value = 3056.4666
$ws.Range("I73").Value = 4050
$ws.Range("J73").Value = 2695.182
$ws.Range("K73").Value = 12150
$ws.Range("L73").Value = 8085.545999999999
$ws.Range("M73").Value = -11214
$ws.Range("N73").Value = -9957.545999999998
# Row 113
$ws.Range("H113").Value = 4912.1665
$ws.Range("I113").Value = 3799.7058
$ws.Range("K113").Value = 3799.7058
$ws.Range("M113").Value = -545.7058000000002
# Row 132
$ws.Range("H132").Value = 3265.7144
$ws.Range("I132").Value = 3745.3333
$ws.Range("J132").Value = 2066.6667
$ws.Range("K132").Value = 11235.9999
$ws.Range("L132").Value = 6200.000100000001
$ws.Range("M132").Value = -8705.999899999999
$ws.Range("N132").Value = -11260.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1798.4736
$ws.Range("I45").Value = 1436.1538
$ws.Range("K45").Value = 1436.1538
$ws.Range("M45").Value = -1059.1538
# Row 61
$ws.Range("H61").Value = 8010.9414
$ws.Range("I61").Value = 8395.691999999999
$ws.Range("J61").Value = 6760.5
$ws.Range("K61").Value = 8395.691999999999
$ws.Range("L61").Value = 6760.5
$ws.Range("M61").Value = -8183.691999999999
$ws.Range("N61").Value = -7184.5
# Row 122
$ws.Range("H122").Value = 4313.0625
$ws.Range("I122").Value = 1851
$ws.Range("J122").Value = 8416.5
$ws.Range("K122").Value = 5553
$ws.Range("L122").Value = 25249.5
$ws.Range("M122").Value = -3103
$ws.Range("N122").Value = -30149.5
# Row 132
$ws.Range("H132").Value = 5519
$ws.Range("I132").Value = 2136.6316
$ws.Range("K132").Value = 6409.8948
$ws.Range("M132").Value = -3879.8948
# Row 136
$ws.Range("H136").Value = 8010.9414
$ws.Range("I136").Value = 8395.691999999999
$ws.Range("J136").Value = 6760.5
$ws.Range("K136").Value = 25187.076
$ws.Range("L136").Value = 20281.5
$ws.Range("M136").Value = -22637.076
$ws.Range("N136").Value = -25381.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7953.0347
$ws.Range("I31").Value = 3402.0527
$ws.Range("K31").Value = 3402.0527
$ws.Range("M31").Value = -3107.0527
# Row 34
$ws.Range("H34").Value = 7953.0347
$ws.Range("I34").Value = 3402.0527
$ws.Range("K34").Value = 3402.0527
$ws.Range("M34").Value = -3200.0527
# Row 41
$ws.Range("H41").Value = 17033.545
$ws.Range("J41").Value = 42499.668
$ws.Range("L41").Value = 42499.668
$ws.Range("N41").Value = -43355.668
# Row 50
$ws.Range("H50").Value = 46084.285
$ws.Range("J50").Value = 67498
$ws.Range("L50").Value = 67498
$ws.Range("N50").Value = -68748
# Row 58
$ws.Range("H58").Value = 3657.8125
$ws.Range("I58").Value = 2705.5789
$ws.Range("J58").Value = 5049.5386
$ws.Range("K58").Value = 2705.5789
$ws.Range("L58").Value = 5049.5386
$ws.Range("M58").Value = -2502.5789
$ws.Range("N58").Value = -5455.5386
# Row 59
$ws.Range("H59").Value = 80399.5
$ws.Range("I59").Value = 89000
$ws.Range("J59").Value = 79443.89
$ws.Range("K59").Value = 89000
$ws.Range("L59").Value = 79443.89
$ws.Range("M59").Value = -87855
$ws.Range("N59").Value = -81733.89
# Row 60
$ws.Range("H60").Value = 13860.277
$ws.Range("J60").Value = 46498.5
$ws.Range("L60").Value = 46498.5
$ws.Range("N60").Value = -47520.5
# Row 94
$ws.Range("H94").Value = 4934.75
$ws.Range("I94").Value = 5185.5
$ws.Range("J94").Value = 4809.375
$ws.Range("K94").Value = 5185.5
$ws.Range("L94").Value = 4809.375
$ws.Range("M94").Value = -4734.5
$ws.Range("N94").Value = -5711.375
# Row 99
$ws.Range("H99").Value = 2483.5
$ws.Range("I99").Value = 2481.6667
$ws.Range("K99").Value = 2481.6667
$ws.Range("M99").Value = -983.6667000000002
# Row 105
$ws.Range("H105").Value = 5577.5
$ws.Range("I105").Value = 4770
$ws.Range("K105").Value = 4770
$ws.Range("M105").Value = -3023
# Row 107
$ws.Range("H107").Value = 1218.0769
$ws.Range("I107").Value = 1003.4286
$ws.Range("J107").Value = 2119.6
$ws.Range("K107").Value = 1003.4286
$ws.Range("L107").Value = 2119.6
$ws.Range("M107").Value = 916.5714
$ws.Range("N107").Value = -5959.6
# Row 122
$ws.Range("H122").Value = 2851
$ws.Range("I122").Value = 2882.9092
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 8648.7276
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6198.7276
$ws.Range("N122").Value = -12400
# Row 126
$ws.Range("H126").Value = 2483.5
$ws.Range("I126").Value = 2481.6667
$ws.Range("K126").Value = 7445.000100000001
$ws.Range("M126").Value = -4975.000100000001
# Row 132
$ws.Range("H132").Value = 4349.125
$ws.Range("I132").Value = 3542.0715
$ws.Range("J132").Value = 9998.5
$ws.Range("K132").Value = 10626.2145
$ws.Range("L132").Value = 29995.5
$ws.Range("M132").Value = -8096.2145
$ws.Range("N132").Value = -35055.5
# Row 136
$ws.Range("H136").Value = 3657.8125
$ws.Range("I136").Value = 2705.5789
$ws.Range("J136").Value = 5049.5386
$ws.Range("K136").Value = 8116.736699999999
$ws.Range("L136").Value = 15148.6158
$ws.Range("M136").Value = -5566.736699999999
$ws.Range("N136").Value = -20248.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 562.6667
$ws.Range("I8").Value = 562.6667
$ws.Range("K8").Value = 1688.0001
$ws.Range("M8").Value = -1549.0001
# Row 117
$ws.Range("H117").Value = 1916.375
$ws.Range("I117").Value = 1500
$ws.Range("J117").Value = 2055.1667
$ws.Range("K117").Value = 4500
$ws.Range("L117").Value = 6165.500100000001
$ws.Range("M117").Value = -1058
$ws.Range("N117").Value = -13049.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 43605
$ws.Range("J26").Value = 43605
$ws.Range("L26").Value = 43605
$ws.Range("N26").Value = -44165
# Row 50
$ws.Range("H50").Value = 43605
$ws.Range("J50").Value = 43605
$ws.Range("L50").Value = 43605
$ws.Range("N50").Value = -44601
# Row 58
$ws.Range("H58").Value = 29000
$ws.Range("I58").Value = 33500
$ws.Range("K58").Value = 33500
$ws.Range("M58").Value = -33223
# Row 122
$ws.Range("H122").Value = 4965.579
$ws.Range("I122").Value = 988.38464
$ws.Range("J122").Value = 13582.833
$ws.Range("K122").Value = 2965.15392
$ws.Range("L122").Value = 40748.499
$ws.Range("M122").Value = -515.1539199999997
$ws.Range("N122").Value = -45648.499

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 5041.778
$ws.Range("I122").Value = 4263.467
$ws.Range("K122").Value = 12790.401
$ws.Range("M122").Value = -10340.401
# Row 136
$ws.Range("H136").Value = 4972.25
$ws.Range("I136").Value = 1899.5834
$ws.Range("K136").Value = 5698.7502
$ws.Range("M136").Value = -3148.7502

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2756
$ws.Range("I122").Value = 2582.5
$ws.Range("J122").Value = 3588.8
$ws.Range("K122").Value = 7747.5
$ws.Range("L122").Value = 10766.4
$ws.Range("M122").Value = -5297.5
$ws.Range("N122").Value = -15666.4
# Row 132
$ws.Range("H132").Value = 2480.5952
$ws.Range("I132").Value = 1904.625
$ws.Range("K132").Value = 5713.875
$ws.Range("M132").Value = -3183.875
# Row 136
$ws.Range("H136").Value = 3700.2307
$ws.Range("I136").Value = 2826.842
$ws.Range("J136").Value = 6070.857
$ws.Range("K136").Value = 8480.526
$ws.Range("L136").Value = 18212.571
$ws.Range("M136").Value = -5930.526
$ws.Range("N136").Value = -23312.571
